$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45947
$ws.Range("B2").Value = 123.53
$ws.Range("C2").Value = 109.99
$ws.Range("D2").Value = 106.86
$ws.Range("E2").Value = 105.16
$ws.Range("F2").Value = 105.05
$ws.Range("G2").Value = 106.83
$ws.Range("H2").Value = 116.32
$ws.Range("I2").Value = 137.28
$ws.Range("J2").Value = 159.4
$ws.Range("K2").Value = 129.77
$ws.Range("L2").Value = 99.38
$ws.Range("M2").Value = 85.22
$ws.Range("N2").Value = 76.79000000000001
$ws.Range("O2").Value = 60.42
$ws.Range("P2").Value = 55.01
$ws.Range("Q2").Value = 55.03
$ws.Range("R2").Value = 68.44
$ws.Range("S2").Value = 94.70999999999999
$ws.Range("T2").Value = 120.75
$ws.Range("U2").Value = 138.19
$ws.Range("V2").Value = 157.65
$ws.Range("W2").Value = 147.76
$ws.Range("X2").Value = 127.88
$ws.Range("Y2").Value = 122.22
$ws.Range("Z2").Value = 108.74
$ws.Range("AB2").Value = 138.88
$ws.Range("AD2").Value = 152.7
$ws.Range("AF2").Value = 144.58
$ws.Range("AG2").Value = "2h-17h"
